$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("B18").Value = "-"
$ws.Range("C18").Value = "[-, -, 'ELM-2NA-Eletropneumática', 'ELM-2NA-Eletro']"
$ws.Range("E18").Value = "['MEC-2NA-C.L.P.', -, 'MEC-2NA-C.pneumática', 'MEC-2NA-C. Hidráulica']"
$ws.Range("F18").Value = "[-, 'MEC-2NB-C. L. P.', -, 'MEC-2NB-Coman. Hidraulicos']"

# Row 19
$ws.Range("B19").Value = "-"
$ws.Range("C19").Value = "[-, -, Guilherme-Eletropneumática-2NA, 'ELM-2NA-Eletro']"
$ws.Range("D19").Value = "ELM-1NA-Máquinas Térmicas e de Fluxo"
$ws.Range("E19").Value = "['MEC-2NA-C.L.P.', -, 'MEC-2NA-C.pneumática', 'MEC-2NA-C. Hidráulica']"
$ws.Range("F19").Value = "[-, 'MEC-2NB-C. L. P.', -, 'MEC-2NB-Coman. Hidraulicos']"

# Row 20
$ws.Range("B20").Value = "-"
$ws.Range("C20").Value = "[-, -, 'ELM-2NA-Eletropneumática', 'ELM-2NA-Eletro']"
$ws.Range("E20").Value = "['MEC-2NA-C.L.P.', -, 'MEC-2NA-C.pneumática', 'MEC-2NA-C. Hidráulica']"
$ws.Range("F20").Value = "[-, 'MEC-2NB-C. L. P.', -, 'MEC-2NB-Coman. Hidraulicos']"

# Row 21
$ws.Range("B21").Value = "-"
$ws.Range("C21").Value = "[-, -, 'ELM-2NA-Eletropneumática', 'ELM-2NA-Eletro']"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "['MEC-2NA-C.L.P.', -, 'MEC-2NA-C.pneumática', 'MEC-2NA-C. Hidráulica']"
$ws.Range("F21").Value = "[-, 'MEC-2NB-C. L. P.', -, 'MEC-2NB-Coman. Hidraulicos']"
